$wb = $excel.ActiveWorkbook

# Original first sheet (currently "Sheet1" with login/task data)
$ws1 = $wb.Worksheets.Item(1)

# Insert a new worksheet before the existing Sheet1; it becomes the active sheet
$ws2 = $wb.Worksheets.Add($ws1)
$ws2.Name = "Sheet2"

# Populate the new sheet with the FTNames list
$ws2.Range("A1").Value = "FTNames"
$ws2.Range("A2").Value = "DecnFT1"
$ws2.Range("A3").Value = "DecnFT2"
$ws2.Range("A4").Value = "DecnFT3"

# Make the header bold
$ws2.Range("A1").Font.Bold = $true

# Ensure the new sheet is the active/selected tab
$ws2.Activate()
